$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.364.69"
$ws.Range("E2").Value = "  -1.40%  "

$ws.Range("D3").Value = "2.993.82"
$ws.Range("E3").Value = "  +0.15%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'506.16"
$ws.Range("E5").Value = "  +1.03%  "

$ws.Range("D6").Value = "'138.01"
$ws.Range("E6").Value = "  -0.09%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "'0.430"
$ws.Range("E8").Value = "  -0.19%  "

$ws.Range("E9").Value = "  -2.19%  "

$ws.Range("E10").Value = "  -0.26%  "

$ws.Range("D11").Value = "'0.366"
$ws.Range("E11").Value = "  +2.00%  "

$ws.Range("D12").Value = "3.503.43"
$ws.Range("E12").Value = "  +0.01%  "

$ws.Range("E13").Value = "  -1.01%  "

$ws.Range("D14").Value = "'25.70"
$ws.Range("E14").Value = "  -1.86%  "

$ws.Range("D15").Value = "'0.0000163"
$ws.Range("E15").Value = "  +1.86%  "

$ws.Range("D16").Value = "56.339.01"
$ws.Range("E16").Value = "  -1.53%  "

$ws.Range("D17").Value = "2.994.31"
$ws.Range("E17").Value = "  -0.19%  "

$ws.Range("D18").Value = "'5.98"
$ws.Range("E18").Value = "  -1.56%  "

$ws.Range("D19").Value = "'12.94"
$ws.Range("E19").Value = "  +2.43%  "

$ws.Range("D20").Value = "'8.06"
$ws.Range("E20").Value = "  +2.19%  "

$ws.Range("D21").Value = "'331.86"
$ws.Range("E21").Value = "  +3.12%  "

$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.12%  "

$ws.Range("D23").Value = "'0.494"
$ws.Range("E23").Value = "  +0.61%  "

$ws.Range("D24").Value = "'64.69"
$ws.Range("E24").Value = "  +2.14%  "

$ws.Range("D25").Value = "3.114.57"
$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").Value = "0.0₃0947"
$ws.Range("E26").Value = "  +5.81%  "

$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  -0.12%  "

$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "'0.163"
$ws.Range("E28").Value = "  +0.49%  "

$ws.Range("D29").Value = "'6.37"
$ws.Range("E29").Value = "  -3.05%  "

$ws.Range("D30").Value = "'6.92"
$ws.Range("E30").Value = "  -2.91%  "

$ws.Range("E31").Value = "  -0.02%  "

$ws.Range("D32").Value = "'20.27"
$ws.Range("E32").Value = "  +0.19%  "

$ws.Range("D33").Value = "'1.16"
$ws.Range("E33").Value = "  -0.80%  "

$ws.Range("D34").Value = "'152.66"
$ws.Range("E34").Value = "  -1.77%  "

$ws.Range("D35").Value = "'4.49"
$ws.Range("E35").Value = "  -1.96%  "

$ws.Range("D36").Value = "'5.80"
$ws.Range("E36").Value = "  +0.32%  "

$ws.Range("D37").Value = "'26.12"
$ws.Range("E37").Value = "  +7.20%  "

$ws.Range("D38").Value = "'1.25"
$ws.Range("E38").Value = "  +0.16%  "

$ws.Range("D39").Value = "'0.0660"
$ws.Range("E39").Value = "  -0.57%  "

$ws.Range("D40").Value = "3.030.56"
$ws.Range("E40").Value = "  +0.34%  "

$ws.Range("E41").Value = "  -2.53%  "

$ws.Range("E42").Value = "  -0.12%  "

$ws.Range("D43").Value = "'3.79"
$ws.Range("E43").Value = "  +0.94%  "

$ws.Range("E44").Value = "  +1.00%  "

$ws.Range("D45").Value = "2.186.28"
$ws.Range("E45").Value = "  -0.46%  "

$ws.Range("E46").Value = "  -2.96%  "

$ws.Range("D47").Value = "'5.84"
$ws.Range("E47").Value = "  -1.64%  "

$ws.Range("D48").Value = "'0.922"
$ws.Range("E48").Value = "  -1.74%  "

$ws.Range("D49").Value = "'0.0235"
$ws.Range("E49").Value = "  -0.15%  "

$ws.Range("D50").Value = "'19.51"
$ws.Range("E50").Value = "  +1.08%  "

$ws.Range("D51").Value = "'0.0850"
$ws.Range("E51").Value = "  -2.59%  "
